$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 227, shifting existing rows 227-286 down to 228-287.
$ws.Rows.Item(227).EntireRow.Insert()

# Populate the newly inserted row 227 with the new record.
$ws.Range("A227").Value = 6
$ws.Range("B227").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C227").Value = "Metropolitana"
$ws.Range("D227").Value = "11/22/2021"
$ws.Range("E227").Value = 13
$ws.Range("F227").Value = 100112030
$ws.Range("G227").Value = "Poroto granado"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 95
$ws.Range("K227").Value = 23000
$ws.Range("L227").Value = 25000
$ws.Range("M227").Value = 24053
$ws.Range("N227").Value = "`$/caja 15 kilos"
$ws.Range("O227").Value = "Provincia de Limarí"
$ws.Range("P227").Value = 1604
$ws.Range("Q227").Value = 15
$ws.Range("R227").Value = "Hortaliza"

Write-Output "done"
